$d = $word.ActiveDocument

# --- 1. "Programa resumido" paragraph: split one run into four via manual line breaks ---
$d.Content.Find.Execute(
    "Análise tridimensional de tensõesCritérios de FalhaFundamentos da Teoria da ElasticidadeAnálise Numérica de Tensões e Deformações",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Análise tridimensional de tensões^lCritérios de Falha^lFundamentos da Teoria da Elasticidade^lAnálise Numérica de Tensões e Deformações",
    2
)

# --- 2. "Programa" paragraph: split the long descriptive run into eight via manual line breaks ---
$d.Content.Find.Execute(
    "Análise tridimensional de tensões:Equações de Transformação no caso Triaxial; Tensões Principais: autovalores e autovetores do tensor de tensão; Invariantes do Estado de Tensão; Tensões Octaédricas.Critérios de Falha:Critérios de Fratura para Materiais Frágeis; Critério de Escoamento da Máxima Tensão Cisalhante (Tresca); Critério da Energia de Distorção (Von Mises); Componentes Hidrostático e Desviador do Estado de Tensão. Fundamentos da Teoria da Elasticidade: Estado de Tensão em um Sólido Contínuo; Estado de deformação: Relações Deformação-Deslocamento (equações cinemáticas), deformação em 3 dimensões e os Invariantes da deformação; Equações Diferenciais de Equilíbrio; Equações de Compatibilidade: Interpretações matemática e física; Princípio de Saint-Venant; Problemas Bidimensionais; Equação de Compatibilidade para o caso bidimensional; Relações Básicas em Coordenadas Polares; Aplicação em Problemas Axissimétricos (tubos de paredes grossas); Análise Numérica de Tensões e Deformações:Diferenças Finitas; Introdução ao Método dos Elementos Finitos; Princípio dos Trabalhos Virtuais, o Problema Unidimensional; Problema Bidimensional; Discretização: Funções de aproximação para elementos triangulares; Emprego de programas computacionais na análise de tensões e deformações pelo Método dos Elementos Finitos",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Análise tridimensional de tensões:^lEquações de Transformação no caso Triaxial; Tensões Principais: autovalores e autovetores do tensor de tensão; Invariantes do Estado de Tensão; Tensões Octaédricas.^lCritérios de Falha:^lCritérios de Fratura para Materiais Frágeis; Critério de Escoamento da Máxima Tensão Cisalhante (Tresca); Critério da Energia de Distorção (Von Mises); Componentes Hidrostático e Desviador do Estado de Tensão. ^lFundamentos da Teoria da Elasticidade: ^lEstado de Tensão em um Sólido Contínuo; Estado de deformação: Relações Deformação-Deslocamento (equações cinemáticas), deformação em 3 dimensões e os Invariantes da deformação; Equações Diferenciais de Equilíbrio; Equações de Compatibilidade: Interpretações matemática e física; Princípio de Saint-Venant; Problemas Bidimensionais; Equação de Compatibilidade para o caso bidimensional; Relações Básicas em Coordenadas Polares; Aplicação em Problemas Axissimétricos (tubos de paredes grossas); ^lAnálise Numérica de Tensões e Deformações:^lDiferenças Finitas; Introdução ao Método dos Elementos Finitos; Princípio dos Trabalhos Virtuais, o Problema Unidimensional; Problema Bidimensional; Discretização: Funções de aproximação para elementos triangulares; Emprego de programas computacionais na análise de tensões e deformações pelo Método dos Elementos Finitos",
    2
)

# --- 3. "Critério:" run: split "NS = NP1+NP2; NP1: ...; NP2: ..." into three via manual line breaks ---
$d.Content.Find.Execute(
    "NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "NS = NP1+NP2; ^lNP1: questões da P1 valendo até 4p. no total; ^lNP2: questões da P2 valendo até 6 p. no total.",
    2
)

# --- 4. "Bibliografia" paragraph: split into eight numbered entries separated by blank lines (double break) ---
$d.Content.Find.Execute(
    "1. J.M. GERE. Mecânica dos Materiais. São Paulo: Pioneira Thomson Learning, 2003, 698p.2. M.H. SADD. Elasticity: Theory, Applications and Numerics. Amsterdam: Elsevier, 2005, 461p.3. R.R. CRAIG,Jr. Mecânica dos Materiais. Rio de Janeiro LTC. 2a Ed., 2003, 552p. 4. A.C. UGURAL, S.K. FENSTER. Advanced Strength and Applied Elasticity. New Jersey: Prentice Hall. 4th Ed., 2003, 544p.5. S.P. TIMOSHENKO, J.N. GOODIER. Teoria da Elasticidade. Rio de janeiro: Guanabara Dois. 3a Ed., 1980, 545p.6. A.R. RAGAB, S.E. BAYOUMI. Engineering Solid Mechanics, Fundamentals and Applications. New York: CRC Press, 1999, 921p. 7. POPOV, E. P. Introdução à Mecânica dos Sólidos, São Paulo: Edgard Blücher, 1978, 552p.8. T.M. ATANACKOVIC, A. GURAN. Theory of Elasticity for Scientists and Engineers. New York: Springer Science+Business, 2000, 374p.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1. J.M. GERE. Mecânica dos Materiais. São Paulo: Pioneira Thomson Learning, 2003, 698p.^l^l2. M.H. SADD. Elasticity: Theory, Applications and Numerics. Amsterdam: Elsevier, 2005, 461p.^l^l3. R.R. CRAIG,Jr. Mecânica dos Materiais. Rio de Janeiro LTC. 2a Ed., 2003, 552p. ^l^l4. A.C. UGURAL, S.K. FENSTER. Advanced Strength and Applied Elasticity. New Jersey: Prentice Hall. 4th Ed., 2003, 544p.^l^l5. S.P. TIMOSHENKO, J.N. GOODIER. Teoria da Elasticidade. Rio de janeiro: Guanabara Dois. 3a Ed., 1980, 545p.^l^l6. A.R. RAGAB, S.E. BAYOUMI. Engineering Solid Mechanics, Fundamentals and Applications. New York: CRC Press, 1999, 921p. ^l^l7. POPOV, E. P. Introdução à Mecânica dos Sólidos, São Paulo: Edgard Blücher, 1978, 552p.^l^l8. T.M. ATANACKOVIC, A. GURAN. Theory of Elasticity for Scientists and Engineers. New York: Springer Science+Business, 2000, 374p.",
    2
)
